# Added base for InnerSpectre
# Appends three new localization rows (new dialogue lines) to the "strings"
# sheet, right after the existing data, mirroring the key/value pattern
# already used on this sheet (column A and column B hold the same literal
# text for these new rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strings")
$ws.Activate()

$ws.Range("A126").Value = "Don't you want something more?"
$ws.Range("B126").Value = "Don't you want something more?"

$ws.Range("A127").Value = "Make up for lost time?"
$ws.Range("B127").Value = "Make up for lost time?"

$ws.Range("A128").Value = "Then prove it to me."
$ws.Range("B128").Value = "Then prove it to me."

# Scroll so the newly-added rows are in view, and select the last edited
# cell, matching the author's end state.
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B128").Select()
